# The author had moved on to editing near the "Background" bullet
# (chapter 2 / extra papers), so Word's "last edit location" marker -
# the hidden _GoBack bookmark - needs to move from after "Main Body"
# to right after "Background".
#
# Bookmarks.Add() re-uses the existing "_GoBack" name, so simply adding
# it at the new location automatically removes it from the old one.

$d = $word.ActiveDocument

# Find the "Background" list-item paragraph by its text (robust to
# paragraph-index drift).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "Background") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $endRng = $target.Range
    $endRng.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
    $endRng.Collapse(0) | Out-Null      # collapse to a single point right after "Background"
    $insertPos = $endRng.End

    # Insert a one-character placeholder so the target position is no
    # longer sitting exactly on the paragraph-end boundary, add the
    # bookmark there (now a safe, interior position), then delete the
    # placeholder again - leaving a proper zero-length bookmark right
    # after the paragraph's text, matching Word's own behaviour.
    $endRng.InsertAfter([string][char]1) | Out-Null

    $bmRng = $d.Range($insertPos, $insertPos)
    $d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

    $placeholder = $d.Range($insertPos, $insertPos + 1)
    $placeholder.Delete() | Out-Null
}
